# Auto-generated update of Yojimbo_Profits market-data snapshot values.
# For each affected Leve row, writes the latest currentAveragePrice* /
# LevePrice* / LeveProfit* figures pulled by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2839.975
$ws.Range("I116").Value = 2432.3235
$ws.Range("J116").Value = 5150
$ws.Range("K116").Value = 2432.3235
$ws.Range("L116").Value = 5150
$ws.Range("M116").Value = 1009.6765
$ws.Range("N116").Value = -12034

$ws.Range("H125").Value = 2677.5
$ws.Range("I125").Value = 3800
$ws.Range("J125").Value = 2303.3333
$ws.Range("K125").Value = 34200
$ws.Range("L125").Value = 20729.9997
$ws.Range("M125").Value = -31740
$ws.Range("N125").Value = -25649.9997

$ws.Range("H137").Value = 3561.6123
$ws.Range("I137").Value = 2990.8708
$ws.Range("J137").Value = 4544.5557
$ws.Range("K137").Value = 8972.6124
$ws.Range("L137").Value = 13633.6671
$ws.Range("M137").Value = -6422.6124
$ws.Range("N137").Value = -18733.6671

$ws.Range("H138").Value = 3579.5469
$ws.Range("I138").Value = 2194.7307
$ws.Range("J138").Value = 4527.0527
$ws.Range("K138").Value = 6584.1921
$ws.Range("L138").Value = 13581.1581
$ws.Range("M138").Value = -1444.1921
$ws.Range("N138").Value = -23861.1581

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1925.4108
$ws.Range("I61").Value = 1061.7858
$ws.Range("J61").Value = 4516.2856
$ws.Range("K61").Value = 1061.7858
$ws.Range("L61").Value = 4516.2856
$ws.Range("M61").Value = -849.7858000000001
$ws.Range("N61").Value = -4940.2856

$ws.Range("H122").Value = 1799.7273
$ws.Range("I122").Value = 1295.4736
$ws.Range("J122").Value = 4993.3335
$ws.Range("K122").Value = 3886.4208
$ws.Range("L122").Value = 14980.0005
$ws.Range("M122").Value = -1436.4208
$ws.Range("N122").Value = -19880.0005

$ws.Range("H132").Value = 2374.9243
$ws.Range("I132").Value = 1876.7091
$ws.Range("J132").Value = 4866
$ws.Range("K132").Value = 5630.1273
$ws.Range("L132").Value = 14598
$ws.Range("M132").Value = -3100.1273
$ws.Range("N132").Value = -19658

$ws.Range("H136").Value = 1925.4108
$ws.Range("I136").Value = 1061.7858
$ws.Range("J136").Value = 4516.2856
$ws.Range("K136").Value = 3185.3574
$ws.Range("L136").Value = 13548.8568
$ws.Range("M136").Value = -635.3574000000003
$ws.Range("N136").Value = -18648.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1033.8813
$ws.Range("I134").Value = 819.06525
$ws.Range("J134").Value = 1794
$ws.Range("K134").Value = 2457.19575
$ws.Range("L134").Value = 5382
$ws.Range("M134").Value = 77.80425000000014
$ws.Range("N134").Value = -10452

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1050
$ws.Range("I16").Value = 1050
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1050
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -763

$ws.Range("H18").Value = 114900
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 114900
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 114900
$ws.Range("N18").Value = -115360

$ws.Range("H99").Value = 1801.7858
$ws.Range("I99").Value = 1340.625
$ws.Range("J99").Value = 2416.6667
$ws.Range("K99").Value = 1340.625
$ws.Range("L99").Value = 2416.6667
$ws.Range("M99").Value = 157.375
$ws.Range("N99").Value = -5412.6667

$ws.Range("H113").Value = 1050
$ws.Range("I113").Value = 1050
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1050
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 1120

$ws.Range("H122").Value = 5019.4287
$ws.Range("I122").Value = 6132.2
$ws.Range("J122").Value = 4007.818
$ws.Range("K122").Value = 18396.6
$ws.Range("L122").Value = 12023.454
$ws.Range("M122").Value = -15946.6
$ws.Range("N122").Value = -16923.454

$ws.Range("H126").Value = 1801.7858
$ws.Range("I126").Value = 1340.625
$ws.Range("J126").Value = 2416.6667
$ws.Range("K126").Value = 4021.875
$ws.Range("L126").Value = 7250.000100000001
$ws.Range("M126").Value = -1551.875
$ws.Range("N126").Value = -12190.0001

$ws.Range("H134").Value = 1732.4193
$ws.Range("I134").Value = 1857.409
$ws.Range("J134").Value = 1426.8889
$ws.Range("K134").Value = 5572.227000000001
$ws.Range("L134").Value = 4280.6667
$ws.Range("M134").Value = -3037.227000000001
$ws.Range("N134").Value = -9350.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 711.44446
$ws.Range("I107").Value = 527.5
$ws.Range("J107").Value = 764
$ws.Range("K107").Value = 1582.5
$ws.Range("L107").Value = 2292
$ws.Range("M107").Value = 337.5
$ws.Range("N107").Value = -6132

$ws.Range("H132").Value = 1372.2222
$ws.Range("I132").Value = 899.44446
$ws.Range("J132").Value = 1845
$ws.Range("K132").Value = 8095.00014
$ws.Range("L132").Value = 16605
$ws.Range("M132").Value = -5565.00014
$ws.Range("N132").Value = -21665

$ws.Range("H139").Value = 1596.2106
$ws.Range("I139").Value = 802.3333
$ws.Range("J139").Value = 2957.1428
$ws.Range("K139").Value = 2406.9999
$ws.Range("L139").Value = 8871.428400000001
$ws.Range("M139").Value = 2733.0001
$ws.Range("N139").Value = -19151.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1411.3
$ws.Range("I102").Value = 1215.0667
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1215.0667
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 406.9332999999999
$ws.Range("N102").Value = -5244

$ws.Range("H113").Value = 1139.1
$ws.Range("I113").Value = 1127.2858
$ws.Range("J113").Value = 1166.6666
$ws.Range("K113").Value = 1127.2858
$ws.Range("L113").Value = 1166.6666
$ws.Range("M113").Value = 1042.7142
$ws.Range("N113").Value = -5506.6666

$ws.Range("H122").Value = 2160.9565
$ws.Range("I122").Value = 1510.1818
$ws.Range("J122").Value = 2757.5
$ws.Range("K122").Value = 4530.5454
$ws.Range("L122").Value = 8272.5
$ws.Range("M122").Value = -2080.5454
$ws.Range("N122").Value = -13172.5

$ws.Range("H126").Value = 1495
$ws.Range("I126").Value = 1490
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 4470
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -2000
$ws.Range("N126").Value = -9440

$ws.Range("H132").Value = 1471.6448
$ws.Range("I132").Value = 1194.695
$ws.Range("J132").Value = 2432.8235
$ws.Range("K132").Value = 3584.085
$ws.Range("L132").Value = 7298.470499999999
$ws.Range("M132").Value = -1054.085
$ws.Range("N132").Value = -12358.4705

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2805.875
$ws.Range("I40").Value = 2528.2222
$ws.Range("J40").Value = 3162.8572
$ws.Range("K40").Value = 2528.2222
$ws.Range("L40").Value = 3162.8572
$ws.Range("M40").Value = -2392.2222
$ws.Range("N40").Value = -3434.8572

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").ClearContents()
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = 0

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 0

$ws.Range("H122").Value = 7777
$ws.Range("I122").Value = 20000
$ws.Range("J122").Value = 4721.25
$ws.Range("K122").Value = 60000
$ws.Range("L122").Value = 14163.75
$ws.Range("M122").Value = -57550
$ws.Range("N122").Value = -19063.75

$ws.Range("H136").Value = 2853.0408
$ws.Range("I136").Value = 2190.3103
$ws.Range("J136").Value = 3814
$ws.Range("K136").Value = 6570.9309
$ws.Range("L136").Value = 11442
$ws.Range("M136").Value = -4020.9309
$ws.Range("N136").Value = -16542

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 40000
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 40000
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 40000
$ws.Range("N57").Value = -41508

$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -830

$ws.Range("H122").Value = 1668734
$ws.Range("I122").Value = 9999999
$ws.Range("J122").Value = 2481
$ws.Range("K122").Value = 29999997
$ws.Range("L122").Value = 7443
$ws.Range("M122").Value = -29997547
$ws.Range("N122").Value = -12343

$ws.Range("H136").Value = 799.1579
$ws.Range("I136").Value = 558.6
$ws.Range("J136").Value = 1261.7693
$ws.Range("K136").Value = 1675.8
$ws.Range("L136").Value = 3785.3079
$ws.Range("M136").Value = 874.1999999999998
$ws.Range("N136").Value = -8885.3079

Write-Output "Updated 33 leve rows across 8 job sheets."